{"js": "// Replace the 25 division-problem answers in the practice table.\n// Each entry is [rowIndex, colIndex, expectedOldText, newText] \u2014 indices\n// are 0-based table grid coordinates (row 0, 4, 8, 12, 16 hold the\n// populated rows; the rows in between are blank spacer rows).\nconst replacements = [\n  [0, 0, \"58\u00f73=19, 1\", \"46\u00f73=15, 1\"],\n  [0, 1, \"73\u00f78=9, 1\", \"99\u00f77=14, 1\"],\n  [0, 2, \"89\u00f79=9, 8\", \"83\u00f77=11, 6\"],\n  [0, 3, \"92\u00f78=11, 4\", \"40\u00f73=13, 1\"],\n  [0, 4, \"71\u00f76=11, 5\", \"64\u00f75=12, 4\"],\n  [4, 0, \"92\u00f77=13, 1\", \"93\u00f77=13, 2\"],\n  [4, 1, \"59\u00f75=11, 4\", \"58\u00f76=9, 4\"],\n  [4, 2, \"40\u00f79=4, 4\", \"87\u00f73=29, 0\"],\n  [4, 3, \"21\u00f76=3, 3\", \"59\u00f72=29, 1\"],\n  [4, 4, \"48\u00f76=8, 0\", \"85\u00f79=9, 4\"],\n  [8, 0, \"46\u00f77=6, 4\", \"98\u00f75=19, 3\"],\n  [8, 1, \"90\u00f72=45, 0\", \"38\u00f72=19, 0\"],\n  [8, 2, \"80\u00f74=20, 0\", \"66\u00f76=11, 0\"],\n  [8, 3, \"51\u00f74=12, 3\", \"30\u00f74=7, 2\"],\n  [8, 4, \"41\u00f75=8, 1\", \"63\u00f74=15, 3\"],\n  [12, 0, \"77\u00f73=25, 2\", \"58\u00f78=7, 2\"],\n  [12, 1, \"26\u00f76=4, 2\", \"86\u00f79=9, 5\"],\n  [12, 2, \"61\u00f75=12, 1\", \"49\u00f75=9, 4\"],\n  [12, 3, \"29\u00f74=7, 1\", \"42\u00f75=8, 2\"],\n  [12, 4, \"30\u00f74=7, 2\", \"58\u00f77=8, 2\"],\n  [16, 0, \"60\u00f73=20, 0\", \"84\u00f79=9, 3\"],\n  [16, 1, \"71\u00f74=17, 3\", \"73\u00f79=8, 1\"],\n  [16, 2, \"78\u00f77=11, 1\", \"68\u00f75=13, 3\"],\n  [16, 3, \"60\u00f75=12, 0\", \"68\u00f73=22, 2\"],\n  [16, 4, \"77\u00f76=12, 5\", \"67\u00f74=16, 3\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Load every target cell's current text first so we can confirm we're\n// editing the right cell before writing (guards against the fact that a\n// couple of the \"new\" strings collide with \"old\" strings used elsewhere\n// in the table, so edits must be addressed by cell position, not by\n// searching for matching text).\nconst cells = replacements.map(([row, col]) => table.getCell(row, col));\ncells.forEach((cell) => cell.load(\"value\"));\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [row, col, oldText, newText] = replacements[i];\n  const cell = cells[i];\n  if (cell.value !== oldText) {\n    throw new Error(\n      `Unexpected text in cell (${row}, ${col}): got ${JSON.stringify(\n        cell.value\n      )}, expected ${JSON.stringify(oldText)}`\n    );\n  }\n  cell.value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 division-problem answers in the practice table.\n# Each entry is (rowIndex, colIndex, expectedOldText, newText) using\n# 1-based Word COM Table.Cell(row, col) coordinates - rows 1, 5, 9, 13, 17\n# hold the populated data (the rows in between are blank spacer rows).\n$replacements = @(\n  @(1, 1, \"58\u00f73=19, 1\", \"46\u00f73=15, 1\"),\n  @(1, 2, \"73\u00f78=9, 1\", \"99\u00f77=14, 1\"),\n  @(1, 3, \"89\u00f79=9, 8\", \"83\u00f77=11, 6\"),\n  @(1, 4, \"92\u00f78=11, 4\", \"40\u00f73=13, 1\"),\n  @(1, 5, \"71\u00f76=11, 5\", \"64\u00f75=12, 4\"),\n  @(5, 1, \"92\u00f77=13, 1\", \"93\u00f77=13, 2\"),\n  @(5, 2, \"59\u00f75=11, 4\", \"58\u00f76=9, 4\"),\n  @(5, 3, \"40\u00f79=4, 4\", \"87\u00f73=29, 0\"),\n  @(5, 4, \"21\u00f76=3, 3\", \"59\u00f72=29, 1\"),\n  @(5, 5, \"48\u00f76=8, 0\", \"85\u00f79=9, 4\"),\n  @(9, 1, \"46\u00f77=6, 4\", \"98\u00f75=19, 3\"),\n  @(9, 2, \"90\u00f72=45, 0\", \"38\u00f72=19, 0\"),\n  @(9, 3, \"80\u00f74=20, 0\", \"66\u00f76=11, 0\"),\n  @(9, 4, \"51\u00f74=12, 3\", \"30\u00f74=7, 2\"),\n  @(9, 5, \"41\u00f75=8, 1\", \"63\u00f74=15, 3\"),\n  @(13, 1, \"77\u00f73=25, 2\", \"58\u00f78=7, 2\"),\n  @(13, 2, \"26\u00f76=4, 2\", \"86\u00f79=9, 5\"),\n  @(13, 3, \"61\u00f75=12, 1\", \"49\u00f75=9, 4\"),\n  @(13, 4, \"29\u00f74=7, 1\", \"42\u00f75=8, 2\"),\n  @(13, 5, \"30\u00f74=7, 2\", \"58\u00f77=8, 2\"),\n  @(17, 1, \"60\u00f73=20, 0\", \"84\u00f79=9, 3\"),\n  @(17, 2, \"71\u00f74=17, 3\", \"73\u00f79=8, 1\"),\n  @(17, 3, \"78\u00f77=11, 1\", \"68\u00f75=13, 3\"),\n  @(17, 4, \"60\u00f75=12, 0\", \"68\u00f73=22, 2\"),\n  @(17, 5, \"77\u00f76=12, 5\", \"67\u00f74=16, 3\")\n)\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\nforeach ($entry in $replacements) {\n    $row = $entry[0]\n    $col = $entry[1]\n    $oldText = $entry[2]\n    $newText = $entry[3]\n\n    $cell = $table.Cell($row, $col)\n    $cellRange = $cell.Range\n    # Cell ranges carry a trailing cell-mark (CR + BEL); trim it off before\n    # comparing so we can confirm we're editing the right cell (addressed\n    # by position, not by text search, since a couple of the \"new\" values\n    # collide with \"old\" values used elsewhere in the table).\n    $current = $cellRange.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $oldText) {\n        throw \"Unexpected text in cell ($row, $col): got '$current', expected '$oldText'\"\n    }\n    $cellRange.Text = $newText\n}\n"}
